# "Use correct version of FoTOMRAEL from develop branch"
#
# Pulls in the FoTOMRAEL value from the develop branch: the
# "Fraction of Technology Outside Modeled Region Affecting Endogenous
# Learning" input (FoTOMRAEL sheet, cell B2) changes from 0.9 to 0.25.

$wb = $excel.ActiveWorkbook

$wsData = $wb.Worksheets.Item("FoTOMRAEL")
$wsData.Range("B2").Value = 0.25

# Reset FoTOMRAEL's lingering selection (was left on B3) back to the
# top-left cell, then leave the "About" sheet as the active/selected tab,
# matching the state the file was saved in on the develop branch.
[void]$wsData.Range("A1").Select()

$wsAbout = $wb.Worksheets.Item("About")
[void]$wsAbout.Select()
